# Variables can now be referenced from sheet name + cell name (e.g. Sheet1!B5),
# demonstrated by adding a new "angle" input on Sheet1 and a cross-sheet
# formula on Sheet2 that reads it. Also tightens the coefficient used for `y`.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New "angle" input row (row 5) on Sheet1.
$ws1.Range("A5").Value = "angle"
$ws1.Range("A5").HorizontalAlignment = -4152   # xlRight, matches the other row labels in column A
$ws1.Range("B5").Value = 0.12
$ws1.Range("D5").HorizontalAlignment = -4152   # keep column D's label style consistent on the new row

# y's formula / description tightened from 2.12345 to 2.1.
$ws1.Range("E2").Value = "2..1*x"
$ws1.Range("F2").Formula = "=x*2.1"

# Sheet2 now references Sheet1's new cell by sheet name + cell name.
$ws2.Range("A1").Formula = "=Sheet1!B5+100"

# Restore/record the selections seen in the edited workbook.
[void]$ws2.Range("B45").Select()
[void]$ws1.Range("E5").Select()
